$d = $word.ActiveDocument

# Locate the paragraph that marks the start of the block to remove: the
# empty paragraph immediately following "LOB1019: Física II (Requisito
# fraco)", and the end of the block: the paragraph containing the
# "© 2020 ..." footer text. Delete the whole run of paragraphs (including
# their paragraph marks) in one shot.
$startPara = $null
$endPara = $null

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text

    if ($t -like "*Ver no Jupiter Salvar em pdf Salvar em docx*") {
        $startPara = $d.Paragraphs.Item($i - 1)
    }
    if ($t -like "*Powered by Jekyll and Github pages*") {
        $endPara = $p
    }
}

$range = $d.Range($startPara.Range.Start, $endPara.Range.End)
$range.Delete()
